$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.090.15"
$ws.Range("E2").Value = "  +6.79%  "
$ws.Range("D3").Value = "3.365.74"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'414.18"
$ws.Range("E5").Value = "  +4.86%  "
$ws.Range("D6").Value = "'111.63"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("D7").Value = "'0.588"
$ws.Range("E7").Value = "  +3.55%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("D10").Value = "'39.56"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Value = "'0.0993"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").Value = "'0.143"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "3.892.35"
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("D14").Value = "'20.04"
$ws.Range("E14").Value = "  +4.77%  "
$ws.Range("D15").Value = "'8.44"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "3.363.59"
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").Value = "'1.05"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "60.883.41"
$ws.Range("E18").Value = "  +6.94%  "
$ws.Range("D19").Value = "'10.73"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'3.39"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "'0.0000111"
$ws.Range("E21").Value = "  +5.37%  "
$ws.Range("D22").Value = "'13.10"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'303.73"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "'75.26"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'28.89"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'7.86"
$ws.Range("E27").Value = "  +8.66%  "
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").Value = "'8.00"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +6.03%  "
$ws.Range("D31").Value = "'0.115"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("D32").Value = "'2.61"
$ws.Range("E32").Value = "  +23.06%  "
$ws.Range("D33").Value = "'11.50"
$ws.Range("E33").Value = "  +4.36%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("D36").Value = "'0.0509"
$ws.Range("E36").Value = "  +5.37%  "
$ws.Range("D37").Value = "'52.24"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").Value = "'3.14"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("D41").Value = "'0.305"
$ws.Range("E41").Value = "  +7.30%  "
$ws.Range("D42").Value = "'137.14"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").Value = "'0.124"
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("D44").Value = "'1.92"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "'3.96"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "'16.96"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'22.75"
$ws.Range("E47").Value = "  +3.43%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'2.25"
$ws.Range("E48").Value = "  +8.58%  "
$ws.Range("D49").Value = "2.185.93"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "'2.42"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  -2.15%  "
